$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2 through 240 all hold the serial date value
# 46060 and need to be bumped to 46061.
$ws.Range("C2:C240").Value = 46061
